# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '27.119.52'
Set-TextValue 'E2' '  +3.36%  '
Set-TextValue 'D3' '1.658.65'
Set-TextValue 'E3' '  +4.13%  '
Set-TextValue 'E4' '  +0.22%  '
Set-TextValue 'D5' '215.50'
Set-TextValue 'E5' '  +1.75%  '
Set-TextValue 'E6' '  +0.92%  '
Set-TextValue 'E7' '  +0.16%  '
Set-TextValue 'E8' '  +2.09%  '
Set-TextValue 'E9' '  +1.50%  '
Set-TextValue 'D10' '19.56'
Set-TextValue 'E10' '  +3.57%  '
Set-TextValue 'E11' '  +1.38%  '
Set-TextValue 'D12' '1.893.53'
Set-TextValue 'E12' '  +4.24%  '
Set-TextValue 'D13' '1.661.93'
Set-TextValue 'E13' '  +4.19%  '
Set-TextValue 'E14' '  +1.72%  '
Set-TextValue 'D15' '0.519'
Set-TextValue 'E15' '  +3.07%  '
Set-TextValue 'E16' '  +2.02%  '
Set-TextValue 'D17' '27.128.60'
Set-TextValue 'E17' '  +3.54%  '
Set-TextValue 'D18' '238.21'
Set-TextValue 'E18' '  +4.15%  '
Set-TextValue 'D19' '7.83'
Set-TextValue 'E19' '  +3.64%  '
Set-TextValue 'E20' '  +1.11%  '
Set-TextValue 'E21' '  +0.09%  '
Set-TextValue 'D22' '4.43'
Set-TextValue 'E22' '  +4.64%  '
Set-TextValue 'D23' '2.26'
Set-TextValue 'E23' '  +5.26%  '
Set-TextValue 'D24' '9.28'
Set-TextValue 'E24' '  +4.34%  '
Set-TextValue 'D25' '145.82'
Set-TextValue 'E25' '  +0.14%  '
Set-TextValue 'E26' '  +0.18%  '
Set-TextValue 'D27' '7.18'
Set-TextValue 'E27' '  +3.08%  '
Set-TextValue 'E28' '  +1.23%  '
Set-TextValue 'D29' '15.81'
Set-TextValue 'E29' '  +3.32%  '
Set-TextValue 'E30' '  +1.04%  '
Set-TextValue 'E31' '  +1.34%  '
Set-TextValue 'D32' '1.537.80'
Set-TextValue 'E32' '  +5.41%  '
Set-TextValue 'E33' '  +2.63%  '
Set-TextValue 'E34' '  +3.30%  '
Set-TextValue 'E35' '  +7.97%  '
Set-TextValue 'E36' '  +0.01%  '
Set-TextValue 'E37' '  +1.57%  '
Set-TextValue 'D38' '0.889'
Set-TextValue 'E38' '  +8.73%  '
Set-TextValue 'D39' '0.0168'
Set-TextValue 'E39' '  +2.90%  '
Set-TextValue 'E40' '  +3.39%  '
Set-TextValue 'E41' '  +0.10%  '
Set-TextValue 'D42' '2.27'
Set-TextValue 'E42' '  +4.58%  '
Set-TextValue 'D43' '66.17'
Set-TextValue 'E43' '  +9.71%  '
Set-TextValue 'D44' '1.800.51'
Set-TextValue 'E44' '  +4.16%  '
Set-TextValue 'E45' '  +2.68%  '
Set-TextValue 'D46' '0.922'
Set-TextValue 'E46' '  -0.69%  '
Set-TextValue 'E47' '  +2.84%  '
Set-TextValue 'E48' '  +0.74%  '
Set-TextValue 'E49' '  +3.23%  '
Set-TextValue 'D50' '0.0505'
Set-TextValue 'E50' '  +1.00%  '
Set-TextValue 'D51' '0.0974'
Set-TextValue 'E51' '  +3.13%  '
